# Fix: when EV included, REF cost is lower than OPT cost.
# Insert a new vehicle row (all-zero placeholder) above the existing data
# so the data set now runs from row 2 (new) through row 12 (shifted).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2, pushing existing data rows down by one.
$ws.Rows.Item(2).Insert() | Out-Null

# Populate the newly inserted row 2.
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = "electricity"
$ws.Cells.Item(2, 3).Value = 0
$ws.Cells.Item(2, 4).Value = "Wh"
$ws.Cells.Item(2, 5).Value = 188
$ws.Cells.Item(2, 6).Value = "Wh/km"
$ws.Cells.Item(2, 7).Value = 0.95
$ws.Cells.Item(2, 8).Value = 6130
$ws.Cells.Item(2, 9).Value = "W"
$ws.Cells.Item(2, 10).Value = 0.9
$ws.Cells.Item(2, 11).Value = 6130
$ws.Cells.Item(2, 12).Value = "W"
$ws.Cells.Item(2, 13).Value = 0

# Match the selection left behind in the saved file.
$ws.Range("C3").Select() | Out-Null
